# feat: add 2022-Q1 data
#
# 1. Create a new worksheet "2022-Q1" (positioned before "总计") containing the
#    per-fund holdings for the new quarter.
# 2. Insert a new top data row into the "总计" (totals) sheet summarizing the
#    2022-Q1 quarter, shifting the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: build the new "2022-Q1" worksheet (fully populated) before moving it
#         into position -- moving a sheet invalidates earlier references to it.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"

# Match the page setup used by the sibling quarterly sheets (0.75in/1in/0.5in
# margins) instead of leaving the generic Excel defaults on the brand-new sheet.
$newSheet.PageSetup.LeftMargin = 0.75 * 72
$newSheet.PageSetup.RightMargin = 0.75 * 72
$newSheet.PageSetup.TopMargin = 1 * 72
$newSheet.PageSetup.BottomMargin = 1 * 72
$newSheet.PageSetup.HeaderMargin = 0.5 * 72
$newSheet.PageSetup.FooterMargin = 0.5 * 72

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newSheet.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# code, name, fund size, total stock position, position ratio, held value (亿元), rank
$fundRows = @(
    @("005505", "前海开源中药研究精选股票A", "11.44", "91.95", "7.89", "0.9026", 5),
    @("005506", "前海开源中药研究精选股票C", "4.62", "91.95", "7.89", "0.3645", 5)
)

for ($r = 0; $r -lt $fundRows.Length; $r++) {
    $rowNum = $r + 2

    $idxCell = $newSheet.Cells.Item($rowNum, 1)
    $idxCell.Value = $r
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $rowData = $fundRows[$r]
    for ($c = 0; $c -lt 6; $c++) {
        $cell = $newSheet.Cells.Item($rowNum, $c + 2)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c]
    }
    $newSheet.Cells.Item($rowNum, 8).Value = $rowData[6]
}

# Now move the fully-populated sheet so it sits right before "总计".
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet.Move($totalSheet)

# ---------------------------------------------------------------------------
# Step 2: update the "总计" worksheet with a new leading row for 2022-Q1,
#         pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$a2 = $totalSheet.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 1.27

# The leading "A" column is a plain 0-based row index that is stored as data
# (not shifted automatically by Rows.Insert), so renumber the rows that moved
# down to keep 0,1,2,3,4,5 sequential again.
for ($i = 3; $i -le 7; $i++) {
    $totalSheet.Cells.Item($i, 1).Value = $i - 2
}
